# Datos_Lab_Fluidos.xlsx edit: normalize Ensayo columns (B:H) for rows 2-15
# by subtracting the row-2 (baseline) value of each column, then restyle
# those cells to match column A's plain/general style. Also extend the
# existing blank formatting from column B into columns C:H for rows 17-31.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B", "C", "D", "E", "F", "G", "H")

# Capture the baseline (row 2) value for every data column before any
# values are overwritten, since row 2 itself is also being rewritten.
$baseline = @{}
foreach ($col in $cols) {
    $baseline[$col] = $ws.Range($col + "2").Value2
}

# Rows 2 through 15 contain the measured data.
for ($row = 2; $row -le 15; $row++) {
    foreach ($col in $cols) {
        $cellAddr = $col + $row
        $current = $ws.Range($cellAddr).Value2
        $ws.Range($cellAddr).Value = $current - $baseline[$col]
    }
}

# Match the formatting (number format / font / style) of column A, which
# uses the plain "Normal" style, across the rewritten B:H cells.
$ws.Range("A2").Copy() | Out-Null
$ws.Range("B2:H15").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Extend the existing (empty) formatted cells in column B down rows
# 17-31 across columns C:H, matching the style already used by B17.
$ws.Range("B17").Copy() | Out-Null
$ws.Range("C17:H31").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
